$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-unused trailing columns (D:F) entirely, shrinking the table
# from A1:F1 down to A1:C1.
$ws.Range("D1:F1").EntireColumn.Delete()

# Update the remaining three headers to the new names used by the
# attendance-export compile.
$ws.Range("A1").Value = "ID Empleado"
$ws.Range("B1").Value = "Nombre"
$ws.Range("C1").Value = "Fecha/Hora"

# Widen column A a bit now that it holds the employee id.
$ws.Range("A1").EntireColumn.ColumnWidth = 17

# Re-apply the autofilter over the new, smaller header range.
$ws.AutoFilterMode = $false
$ws.Range("A1:C1").AutoFilter()

# The _FilterDatabase defined name needs to follow the autofilter range too.
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Asistencia!`$A`$1:`$C`$1"
    }
}

# Match the saved selection state from the authored workbook.
$ws.Range("C2").Select()
